# Apply the ValueSet-fr-current-medication-document-type.xlsx update:
#  - URL property: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication
#  - Title property: "InterOp'Santé" -> "Interop'Santé"
#  - Date property: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
#  - Jurisdiction property value: "" -> "FRANCE"
#  - System URI (CodeSystem) on the Include sheet: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-current-medication-document-type"
$wsMeta.Range("B5").Value = "value set Interop'Santé - type de document de la ressource Composition d'une FCT"
$wsMeta.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$wsMeta.Range("B11").Value = "FRANCE"

$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-document-type"
